$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 6.48 = 25667.62 pesos
✅ 25667.62 pesos = 6.46 = 962.43 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%"

# --- tasas: update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 154.28
$ws2.Range("O10").Value = 3960
$ws2.Range("N12").Value = 3973.9
$ws2.Range("O12").Value = 149.005
